# Apply scheduled market-data refresh to Behemoth_Profits workbook
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# on each job class sheet with freshly pulled market-board values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6511.4287
$ws.Range("I18").Value = 6511.4287
$ws.Range("K18").Value = 6511.4287
$ws.Range("M18").Value = -6227.4287

$ws.Range("H28").Value = 363.4
$ws.Range("J28").Value = 637
$ws.Range("L28").Value = 637
$ws.Range("N28").Value = -1607

$ws.Range("H43").Value = 7534
$ws.Range("I43").Value = 2124.25
$ws.Range("J43").Value = 9938.333000000001
$ws.Range("K43").Value = 2124.25
$ws.Range("L43").Value = 9938.333000000001
$ws.Range("M43").Value = -2055.25
$ws.Range("N43").Value = -10076.333

$ws.Range("H74").Value = 9999.666999999999
$ws.Range("I74").Value = 7666
$ws.Range("J74").Value = 12333.333
$ws.Range("K74").Value = 7666
$ws.Range("L74").Value = 12333.333
$ws.Range("M74").Value = -6730
$ws.Range("N74").Value = -14205.333

$ws.Range("H77").Value = 9999.666999999999
$ws.Range("I77").Value = 7666
$ws.Range("J77").Value = 12333.333
$ws.Range("K77").Value = 38330
$ws.Range("L77").Value = 61666.665
$ws.Range("M77").Value = -33650
$ws.Range("N77").Value = -71026.66500000001

$ws.Range("H112").Value = 2337.6428
$ws.Range("I112").Value = 1550
$ws.Range("J112").Value = 2468.9167
$ws.Range("K112").Value = 4650
$ws.Range("L112").Value = 7406.750100000001
$ws.Range("M112").Value = -3542
$ws.Range("N112").Value = -9622.750100000001

$ws.Range("H121").Value = 1083789.1
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1083789.1
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 3251367.3
$ws.Range("M121").Value = $null
$ws.Range("N121").Value = -3254861.3

$ws.Range("H129").Value = 860.2
$ws.Range("I129").Value = 575.25
$ws.Range("K129").Value = 1725.75
$ws.Range("M129").Value = 3274.25

$ws.Range("H137").Value = 3983
$ws.Range("I137").Value = 967.0909
$ws.Range("J137").Value = 5934.4707
$ws.Range("K137").Value = 2901.2727
$ws.Range("L137").Value = 17803.4121
$ws.Range("M137").Value = -351.2727
$ws.Range("N137").Value = -22903.4121

$ws.Range("H138").Value = 3262.75
$ws.Range("I138").Value = 3286.25
$ws.Range("J138").Value = 3258.05
$ws.Range("K138").Value = 9858.75
$ws.Range("L138").Value = 9774.150000000001
$ws.Range("M138").Value = -4718.75
$ws.Range("N138").Value = -20054.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7949899.5
$ws.Range("I32").Value = 8481587
$ws.Range("J32").Value = 107503
$ws.Range("K32").Value = 8481587
$ws.Range("L32").Value = 107503
$ws.Range("M32").Value = -8481300
$ws.Range("N32").Value = -108077

$ws.Range("H61").Value = 19280166
$ws.Range("I61").Value = 45460696
$ws.Range("J61").Value = 81109.87
$ws.Range("K61").Value = 45460696
$ws.Range("L61").Value = 81109.87
$ws.Range("M61").Value = -45460484
$ws.Range("N61").Value = -81533.87

$ws.Range("H63").Value = 6072.3335
$ws.Range("I63").Value = 3035.2
$ws.Range("K63").Value = 3035.2
$ws.Range("M63").Value = -2349.2

$ws.Range("H66").Value = 6072.3335
$ws.Range("I66").Value = 3035.2
$ws.Range("K66").Value = 15176
$ws.Range("M66").Value = -11744

$ws.Range("H74").Value = 6415192.5
$ws.Range("I74").Value = 7813506
$ws.Range("J74").Value = 22901.428
$ws.Range("K74").Value = 7813506
$ws.Range("L74").Value = 22901.428
$ws.Range("M74").Value = -7812632
$ws.Range("N74").Value = -24649.428

$ws.Range("H77").Value = 6415192.5
$ws.Range("I77").Value = 7813506
$ws.Range("J77").Value = 22901.428
$ws.Range("K77").Value = 39067530
$ws.Range("L77").Value = 114507.14
$ws.Range("M77").Value = -39063162
$ws.Range("N77").Value = -123243.14

$ws.Range("H110").Value = 2396.6667
$ws.Range("I110").Value = 2276.4
$ws.Range("K110").Value = 2276.4
$ws.Range("M110").Value = -231.4000000000001

$ws.Range("H136").Value = 19280166
$ws.Range("I136").Value = 45460696
$ws.Range("J136").Value = 81109.87
$ws.Range("K136").Value = 136382088
$ws.Range("L136").Value = 243329.61
$ws.Range("M136").Value = -136379538
$ws.Range("N136").Value = -248429.61

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1395.0714
$ws.Range("I20").Value = 1108.7778
$ws.Range("J20").Value = 1910.4
$ws.Range("K20").Value = 1108.7778
$ws.Range("L20").Value = 1910.4
$ws.Range("M20").Value = -861.7778000000001
$ws.Range("N20").Value = -2404.4

$ws.Range("H86").Value = 1808.8948
$ws.Range("I86").Value = 1588.8462
$ws.Range("J86").Value = 2285.6667
$ws.Range("K86").Value = 1588.8462
$ws.Range("L86").Value = 2285.6667
$ws.Range("M86").Value = -465.8462
$ws.Range("N86").Value = -4531.6667

$ws.Range("H89").Value = 1808.8948
$ws.Range("I89").Value = 1588.8462
$ws.Range("J89").Value = 2285.6667
$ws.Range("K89").Value = 7944.231
$ws.Range("L89").Value = 11428.3335
$ws.Range("M89").Value = -2328.231
$ws.Range("N89").Value = -22660.3335

$ws.Range("H105").Value = 1942.1111
$ws.Range("I105").Value = 1355.5714
$ws.Range("J105").Value = 3995
$ws.Range("K105").Value = 1355.5714
$ws.Range("L105").Value = 3995
$ws.Range("M105").Value = 391.4286
$ws.Range("N105").Value = -7489

$ws.Range("H107").Value = 1237.2106
$ws.Range("I107").Value = 1154.7059
$ws.Range("J107").Value = 1938.5
$ws.Range("K107").Value = 1154.7059
$ws.Range("L107").Value = 1938.5
$ws.Range("M107").Value = 765.2941000000001
$ws.Range("N107").Value = -5778.5

$ws.Range("H134").Value = 57321.473
$ws.Range("I134").Value = 4439.6
$ws.Range("J134").Value = 255628.5
$ws.Range("K134").Value = 13318.8
$ws.Range("L134").Value = 766885.5
$ws.Range("M134").Value = -10783.8
$ws.Range("N134").Value = -771955.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1748.625
$ws.Range("I22").Value = 569.8570999999999
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 569.8570999999999
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -219.8570999999999
$ws.Range("N22").Value = -10700

$ws.Range("H31").Value = 873881
$ws.Range("I31").Value = 31355
$ws.Range("J31").Value = 1065364.1
$ws.Range("K31").Value = 31355
$ws.Range("L31").Value = 1065364.1
$ws.Range("M31").Value = -31060
$ws.Range("N31").Value = -1065954.1

$ws.Range("H34").Value = 873881
$ws.Range("I34").Value = 31355
$ws.Range("J34").Value = 1065364.1
$ws.Range("K34").Value = 31355
$ws.Range("L34").Value = 1065364.1
$ws.Range("M34").Value = -31153
$ws.Range("N34").Value = -1065768.1

$ws.Range("H39").Value = 12978.429
$ws.Range("I39").Value = 3616.3333
$ws.Range("K39").Value = 3616.3333
$ws.Range("M39").Value = -3225.3333

$ws.Range("H49").Value = 12978.429
$ws.Range("I49").Value = 3616.3333
$ws.Range("K49").Value = 3616.3333
$ws.Range("M49").Value = -3434.3333

$ws.Range("H135").Value = 117499.5
$ws.Range("J135").Value = 117499.5
$ws.Range("L135").Value = 117499.5
$ws.Range("N135").Value = -127639.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1828.7778
$ws.Range("J68").Value = 1777.3043
$ws.Range("L68").Value = 5331.9129
$ws.Range("N68").Value = -6953.9129

$ws.Range("H71").Value = 1828.7778
$ws.Range("J71").Value = 1777.3043
$ws.Range("L71").Value = 15995.7387
$ws.Range("N71").Value = -24107.7387

$ws.Range("H107").Value = 718.4545000000001
$ws.Range("J107").Value = 1012.0833
$ws.Range("L107").Value = 3036.2499
$ws.Range("N107").Value = -6876.2499

$ws.Range("H113").Value = 2049.0715
$ws.Range("J113").Value = 2133.25
$ws.Range("L113").Value = 6399.75
$ws.Range("N113").Value = -10739.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").Value = $null

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").Value = $null

$ws.Range("H46").Value = 14500
$ws.Range("J46").Value = 14500
$ws.Range("L46").Value = 14500
$ws.Range("N46").Value = -14812

$ws.Range("H69").Value = 79641
$ws.Range("J69").Value = 79641
$ws.Range("L69").Value = 79641
$ws.Range("N69").Value = -81139

$ws.Range("H72").Value = 79641
$ws.Range("J72").Value = 79641
$ws.Range("L72").Value = 238923
$ws.Range("N72").Value = -246411

$ws.Range("H113").Value = 2907.0435
$ws.Range("I113").Value = 2270.9333
$ws.Range("J113").Value = 4099.75
$ws.Range("K113").Value = 2270.9333
$ws.Range("L113").Value = 4099.75
$ws.Range("M113").Value = -100.9333000000001
$ws.Range("N113").Value = -8439.75

$ws.Range("H132").Value = 71431140
$ws.Range("I132").Value = 86959060
$ws.Range("J132").Value = 2674.2
$ws.Range("K132").Value = 260877180
$ws.Range("L132").Value = 8022.599999999999
$ws.Range("M132").Value = -260874650
$ws.Range("N132").Value = -13082.6

$ws.Range("H136").Value = 326775.34
$ws.Range("J136").Value = 326775.34
$ws.Range("L136").Value = 980326.02
$ws.Range("N136").Value = -985426.02

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1382.3125
$ws.Range("I16").Value = 1373.1666
$ws.Range("J16").Value = 1409.75
$ws.Range("K16").Value = 1373.1666
$ws.Range("L16").Value = 1409.75
$ws.Range("M16").Value = -1203.1666
$ws.Range("N16").Value = -1749.75

$ws.Range("H46").Value = 2203.4285
$ws.Range("I46").Value = 2255.2856
$ws.Range("J46").Value = 2151.5715
$ws.Range("K46").Value = 2255.2856
$ws.Range("L46").Value = 2151.5715
$ws.Range("M46").Value = -2067.2856
$ws.Range("N46").Value = -2527.5715

$ws.Range("H93").Value = 71429540
$ws.Range("I93").Value = 76923890
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 76923890
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -76922642
$ws.Range("N93").Value = -5496

$ws.Range("H122").Value = 7049.5454
$ws.Range("I122").Value = 6145.4614
$ws.Range("K122").Value = 18436.3842
$ws.Range("M122").Value = -15986.3842

$ws.Range("H132").Value = 7941049.5
$ws.Range("I132").Value = 1007339.06
$ws.Range("J132").Value = 33364654
$ws.Range("K132").Value = 3022017.18
$ws.Range("L132").Value = 100093962
$ws.Range("M132").Value = -3019487.18
$ws.Range("N132").Value = -100099022

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4009.75
$ws.Range("I136").Value = 1415.6
$ws.Range("K136").Value = 4246.799999999999
$ws.Range("M136").Value = -1696.799999999999
